$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 2.88
$ws.Range("L2").Value = 5.5
$ws.Range("Q2").Value = 2.88
$ws.Range("R2").Value = 1.4
$ws.Range("U2").Value = 2.5
$ws.Range("V2").Value = 1.5
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 7.5
$ws.Range("Y2").Value = 10
$ws.Range("Z2").Value = 17
$ws.Range("AA2").Value = 21
$ws.Range("AD2").Value = 6.5
$ws.Range("AE2").Value = 23
$ws.Range("AH2").Value = 8.5
$ws.Range("AI2").Value = 21
$ws.Range("AJ2").Value = 17
$ws.Range("AK2").Value = 51
$ws.Range("AN2").Value = 3.75
$ws.Range("AO2").Value = 12
$ws.Range("AP2").Value = 29
$ws.Range("AQ2").Value = 41
$ws.Range("AR2").Value = 81
$ws.Range("AW2").Value = 6
$ws.Range("AX2").Value = 29
$ws.Range("AZ2").Value = 126
$ws.Range("BA2").Value = 201
$ws.Range("BB2").Value = 501
$ws.Range("G3").Value = 2.88
$ws.Range("I3").Value = 2.88
$ws.Range("L3").Value = 4
$ws.Range("AA3").Value = 34
$ws.Range("AI3").Value = 12
$ws.Range("AK3").Value = 34
$ws.Range("AW3").Value = 4.5
$ws.Range("AZ3").Value = 81
$ws.Range("L7").Value = 7.5
$ws.Range("Q7").Value = 1.68
$ws.Range("R7").Value = 2.03
$ws.Range("Q8").Value = 1.62
$ws.Range("R8").Value = 2.1
$ws.Range("G9").Value = 2.9
$ws.Range("I9").Value = 2.55
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("Q9").Value = 2.15
$ws.Range("R9").Value = 1.63
$ws.Range("W9").Value = 8
$ws.Range("AK9").Value = 26
$ws.Range("AL9").Value = 23
$ws.Range("AV9").Value = 51
$ws.Range("G10").Value = 2.05
$ws.Range("H10").Value = 3.1
$ws.Range("I10").Value = 4
$ws.Range("J10").Value = 2.88
$ws.Range("L10").Value = 5
$ws.Range("M10").Value = 1.14
$ws.Range("N10").Value = 5.5
$ws.Range("Q10").Value = 2.87
$ws.Range("R10").Value = 1.37
$ws.Range("X10").Value = 8
$ws.Range("Y10").Value = 10
$ws.Range("AN10").Value = 3.75
$ws.Range("AP10").Value = 29
$ws.Range("AQ10").Value = 41
$ws.Range("G11").Value = 1.62
$ws.Range("H11").Value = 3.6
$ws.Range("I11").Value = 6
$ws.Range("J11").Value = 2.3
$ws.Range("L11").Value = 6.5
$ws.Range("U11").Value = 2.38
$ws.Range("V11").Value = 1.53
$ws.Range("Z11").Value = 11
$ws.Range("AE11").Value = 23
$ws.Range("AI11").Value = 29
$ws.Range("AJ11").Value = 21
$ws.Range("AO11").Value = 8.5
$ws.Range("AQ11").Value = 29
$ws.Range("AR11").Value = 51
$ws.Range("AZ11").Value = 151
$ws.Range("G12").Value = 2.25
$ws.Range("I12").Value = 3.4
$ws.Range("J12").Value = 3.1
$ws.Range("O12").Value = 1.44
$ws.Range("P12").Value = 2.63
$ws.Range("AI12").Value = 15
$ws.Range("AQ12").Value = 51
$ws.Range("AZ12").Value = 67
$ws.Range("BA12").Value = 101
$ws.Range("BB12").Value = 301
$ws.Range("Q14").Value = 1.75
$ws.Range("R14").Value = 2.05
